# Add the newest daily ORA-errors record (2025-12-04, count 12) as row 28,
# directly below the existing last row (27: 2025-12-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 27's formatting (the date-formatted style) down into row 28 first,
# so the new date cell picks up the same number format/style index as the
# rest of column A instead of Excel minting a brand-new style entry.
$ws.Range("A27").Copy($ws.Range("A28")) | Out-Null

# Date serial 45995 == 12/4/2025 (one day after the previous last row, 45994).
$ws.Range("A28").Value = 45995
$ws.Range("B28").Value = 12

# Match the author's final selection/active cell on the newly added row.
$ws.Range("A28:B28").Select() | Out-Null
